$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 15498.625
$ws.Range("I51").Value = 12997.25
$ws.Range("K51").Value = 12997.25
$ws.Range("M51").Value = -12513.25

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13333
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("H95").Value = 35208
$ws.Range("J95").Value = 35208
$ws.Range("L95").Value = 35208
$ws.Range("N95").Value = -40700
$ws.Range("H97").Value = 66670772
$ws.Range("I97").Value = 83337220
$ws.Range("K97").Value = 83337220
$ws.Range("M97").Value = -83336724
$ws.Range("H102").Value = 26251300
$ws.Range("I102").Value = 1667733.4
$ws.Range("K102").Value = 1667733.4
$ws.Range("M102").Value = -1666111.4
$ws.Range("H122").Value = 4080.5454
$ws.Range("I122").Value = 4098.6
$ws.Range("K122").Value = 12295.8
$ws.Range("M122").Value = -9845.800000000001
$ws.Range("H132").Value = 2557.6667
$ws.Range("I132").Value = 2557.6667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7673.000100000001
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -5143.000100000001
$ws.Range("M32").ClearContents()
$ws.Range("M132").ClearContents()

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 221020.6
$ws.Range("I94").Value = 367769.34
$ws.Range("J94").Value = 897.5
$ws.Range("K94").Value = 367769.34
$ws.Range("L94").Value = 897.5
$ws.Range("M94").Value = -367318.34
$ws.Range("N94").Value = -1799.5
$ws.Range("H99").Value = 2178.5
$ws.Range("I99").Value = 2178.5
$ws.Range("K99").Value = 2178.5
$ws.Range("M99").Value = -680.5
$ws.Range("H107").Value = 100757.5
$ws.Range("I107").Value = 100757.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 100757.5
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = -98837.5
$ws.Range("M107").ClearContents()

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("H16").Value = 399.07693
$ws.Range("I16").Value = 288.8
$ws.Range("J16").Value = 766.6667
$ws.Range("K16").Value = 288.8
$ws.Range("L16").Value = 766.6667
$ws.Range("M16").Value = -1.800000000000011
$ws.Range("N16").Value = -1340.6667
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("N64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("N67").Value = 0
$ws.Range("H107").Value = 837
$ws.Range("I107").Value = 453.66666
$ws.Range("J107").Value = 1124.5
$ws.Range("K107").Value = 453.66666
$ws.Range("L107").Value = 1124.5
$ws.Range("M107").Value = 1466.33334
$ws.Range("N107").Value = -4964.5
$ws.Range("H113").Value = 399.07693
$ws.Range("I113").Value = 288.8
$ws.Range("J113").Value = 766.6667
$ws.Range("K113").Value = 288.8
$ws.Range("L113").Value = 766.6667
$ws.Range("M113").Value = 1881.2
$ws.Range("N113").Value = -5106.6667
$ws.Range("H134").Value = 3445.1428
$ws.Range("I134").Value = 1861.0588
$ws.Range("K134").Value = 5583.1764
$ws.Range("M134").Value = -3048.1764
$ws.Range("H141").Value = 503653.3
$ws.Range("J141").Value = 503653.3
$ws.Range("L141").Value = 503653.3
$ws.Range("N141").Value = -514013.3
$ws.Range("L9").ClearContents()
$ws.Range("L64").ClearContents()
$ws.Range("L67").ClearContents()

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1796.625
$ws.Range("I14").Value = 1796.625
$ws.Range("K14").Value = 5389.875
$ws.Range("M14").Value = -5216.875
$ws.Range("H23").Value = 44.666668
$ws.Range("J23").Value = 42.5
$ws.Range("L23").Value = 127.5
$ws.Range("N23").Value = -597.5
$ws.Range("H38").Value = 322.14285
$ws.Range("J38").Value = 681.6667
$ws.Range("L38").Value = 2045.0001
$ws.Range("N38").Value = -2739.0001
$ws.Range("H107").Value = 959.1667
$ws.Range("I107").Value = 483.66666
$ws.Range("J107").Value = 1434.6666
$ws.Range("K107").Value = 1450.99998
$ws.Range("L107").Value = 4303.9998
$ws.Range("M107").Value = 469.0000199999999
$ws.Range("N107").Value = -8143.9998

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1904.6
$ws.Range("I132").Value = 1904.6
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5713.799999999999
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -3183.799999999999
$ws.Range("M132").ClearContents()

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("H7").Value = 1086.4286
$ws.Range("I7").Value = 1086.4286
$ws.Range("K7").Value = 1086.4286
$ws.Range("M7").Value = -974.4286
$ws.Range("H16").Value = 1466.2
$ws.Range("J16").Value = 1424
$ws.Range("L16").Value = 1424
$ws.Range("N16").Value = -1764
$ws.Range("H20").Value = 6145.5713
$ws.Range("I20").Value = 6145.5713
$ws.Range("K20").Value = 6145.5713
$ws.Range("M20").Value = -5919.5713
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("H126").Value = 1086.4286
$ws.Range("I126").Value = 1086.4286
$ws.Range("K126").Value = 3259.2858
$ws.Range("M126").Value = -789.2857999999997
$ws.Range("H136").Value = 3995
$ws.Range("I136").Value = 3995
$ws.Range("K136").Value = 11985
$ws.Range("M136").Value = -9435
$ws.Range("L5").ClearContents()
$ws.Range("L22").ClearContents()
$ws.Range("L27").ClearContents()

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("H14").Value = 1002.5
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 1005
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 1005
$ws.Range("M14").Value = -832
$ws.Range("N14").Value = -1341
$ws.Range("H20").Value = 10
$ws.Range("I20").Value = 10
$ws.Range("K20").Value = 10
$ws.Range("M20").Value = 230
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("N24").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("H81").Value = 3352.625
$ws.Range("I81").Value = 3118.7144
$ws.Range("J81").Value = 4990
$ws.Range("K81").Value = 6237.4288
$ws.Range("L81").Value = 9980
$ws.Range("M81").Value = -5176.4288
$ws.Range("N81").Value = -12102
$ws.Range("H84").Value = 3352.625
$ws.Range("I84").Value = 3118.7144
$ws.Range("J84").Value = 4990
$ws.Range("K84").Value = 31187.144
$ws.Range("L84").Value = 49900
$ws.Range("M84").Value = -25883.144
$ws.Range("N84").Value = -60508
$ws.Range("H96").Value = 4216.5
$ws.Range("J96").Value = 3933
$ws.Range("L96").Value = 3933
$ws.Range("N96").Value = -6679
$ws.Range("H122").Value = 1601.1875
$ws.Range("J122").Value = 2004.5
$ws.Range("L122").Value = 6013.5
$ws.Range("N122").Value = -10913.5
$ws.Range("H136").Value = 1893.2222
$ws.Range("I136").Value = 1893.2222
$ws.Range("K136").Value = 5679.6666
$ws.Range("M136").Value = -3129.6666
$ws.Range("M11").ClearContents()
$ws.Range("L24").ClearContents()
$ws.Range("M24").ClearContents()
$ws.Range("L30").ClearContents()
